# Apply the "Add files via upload" revision to fins.xlsx.
#
# Summary of the change:
#  - Sheet "6_" has its old "P/A ratio" question content replaced with a new
#    "fin efficiency" question, and becomes the active/selected tab.
#  - Sheet "5_" (previously blank) gets a new "fin effectiveness" question.
#  - Sheet "3_" is no longer the active tab; its stored selection moves to A2.
#  - The shared strings belonging only to the old "6_" content fall out of
#    sharedStrings.xml on save because nothing references them any more.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "6_" : replace the old "P/A ratio" question with the new
# fin-efficiency question.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("6_")

# Wipe the old rows completely (not just their contents) so no stray
# formatted-but-empty cells are left behind from the previous layout.
$ws6.Range("A1:C7").EntireRow.Delete()

$ws6.Range("A1").Value = "Find the fin efficiency of the adiabatic BC fin in our simulation."
$ws6.Range("B1").Value = "Leeway"
$ws6.Range("C1").Value = "Comments"

$ws6.Range("A2").Value = 0.206
$ws6.Range("B2").Value = 0.005

$ws6.Range("C3").Value = "Remember to find the actual heat transfer through the base, not just the flux"
$ws6.Range("C4").Value = "Use the convection rate equation to find the HT through the fin if the fin was a maximum temperature"
$ws6.Range("C5").Value = "Answer should be between 0 and 1.0"
$ws6.Range("C6").Value = "Be sure to find the area of the full fin to get the ideal heat transfer"

# Trailing "leeway" placeholder cells (blank, centered+wrapped style) below
# the question, matching the convention used on the other question sheets.
$ws6.Cells.Item(5, 2).HorizontalAlignment = -4108
$ws6.Cells.Item(5, 2).VerticalAlignment = -4108
$ws6.Cells.Item(5, 2).WrapText = $true

$ws6.Cells.Item(6, 2).HorizontalAlignment = -4108
$ws6.Cells.Item(6, 2).VerticalAlignment = -4108
$ws6.Cells.Item(6, 2).WrapText = $true

$ws6.Cells.Item(7, 2).HorizontalAlignment = -4108
$ws6.Cells.Item(7, 2).VerticalAlignment = -4108
$ws6.Cells.Item(7, 2).WrapText = $true

$ws6.Rows.Item(1).RowHeight = 30
$ws6.Rows.Item(2).RowHeight = 15
$ws6.Rows.Item(3).RowHeight = 45
$ws6.Rows.Item(4).RowHeight = 60
$ws6.Rows.Item(5).RowHeight = 30
$ws6.Rows.Item(6).RowHeight = 30

# ---------------------------------------------------------------------
# Sheet "5_" : add the fin-effectiveness question (previously blank sheet)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("5_")

$ws5.Range("A1").Value = "Find the fin effectiveness of the adiabatic BC fin in our simulation."
$ws5.Range("B1").Value = "Leeway"
$ws5.Range("C1").Value = "Comments"

$ws5.Range("A2").Value = 8.464
$ws5.Range("B2").Value = 0.1

$ws5.Range("C3").Value = "Remember to find the actual heat transfer through the base, not just the flux"
$ws5.Range("C4").Value = "Use the convection rate equation to find the HT through the base if the fin was not there."
$ws5.Range("C5").Value = "Answer should be above 2.0, but likely not above 10.0 or so."

$ws5.Rows.Item(1).RowHeight = 30
$ws5.Rows.Item(3).RowHeight = 45
$ws5.Rows.Item(4).RowHeight = 45
$ws5.Rows.Item(5).RowHeight = 30

# ---------------------------------------------------------------------
# Sheet "3_" : no longer the active tab; selection parked at A2.
# (Select this *before* selecting/activating "6_" below, so "6_" ends up
# as the final active/selected tab -- matching the new activeTab.)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("3_")
$ws3.Range("A2").Select()

# ---------------------------------------------------------------------
# Final selections: "5_" parks on C6, "6_" becomes the active tab on C7.
# ---------------------------------------------------------------------
$ws5.Range("C6").Select()

$ws6.Activate()
$ws6.Range("C7").Select()
